$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 140, pushing the existing
# date-blocks (rows 140:157) down to (142:159).
$ws.Range("A140:A141").EntireRow.Insert()

# Fill the newly inserted row 140 (Calidad = Primera, Fecha = 2022-07-27)
$ws.Cells.Item(140, 1).Value = 11
$ws.Cells.Item(140, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(140, 3).Value = "Bíobío"
$ws.Cells.Item(140, 4).Value = 44769
$ws.Cells.Item(140, 5).Value = 8
$ws.Cells.Item(140, 6).Value = 100112044
$ws.Cells.Item(140, 7).Value = "Perejil"
$ws.Cells.Item(140, 8).Value = "Sin especificar"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 200
$ws.Cells.Item(140, 11).Value = 700
$ws.Cells.Item(140, 12).Value = 800
$ws.Cells.Item(140, 13).Value = 750
$ws.Cells.Item(140, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(140, 15).Value = "Región de Ñuble"
$ws.Cells.Item(140, 16).Value = 750
$ws.Cells.Item(140, 17).Value = 1
$ws.Cells.Item(140, 18).Value = "Hortaliza"

# Fill the newly inserted row 141 (Calidad = Segunda, Fecha = 2022-07-27)
$ws.Cells.Item(141, 1).Value = 11
$ws.Cells.Item(141, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(141, 3).Value = "Bíobío"
$ws.Cells.Item(141, 4).Value = 44769
$ws.Cells.Item(141, 5).Value = 8
$ws.Cells.Item(141, 6).Value = 100112044
$ws.Cells.Item(141, 7).Value = "Perejil"
$ws.Cells.Item(141, 8).Value = "Sin especificar"
$ws.Cells.Item(141, 9).Value = "Segunda"
$ws.Cells.Item(141, 10).Value = 100
$ws.Cells.Item(141, 11).Value = 600
$ws.Cells.Item(141, 12).Value = 600
$ws.Cells.Item(141, 13).Value = 600
$ws.Cells.Item(141, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(141, 15).Value = "Región de Ñuble"
$ws.Cells.Item(141, 16).Value = 600
$ws.Cells.Item(141, 17).Value = 1
$ws.Cells.Item(141, 18).Value = "Hortaliza"

# Make sure the Fecha column keeps the date/time number format used by the
# rest of the column (style index 2 in the original workbook).
$ws.Range("D140:D141").NumberFormat = $ws.Range("D142").NumberFormat
